# Apply updated cryptos data (prices / volume deltas / two swapped rows)
# Values that look like plain numbers are written with a leading apostrophe
# so Excel stores them as text (matching the original inline-string cells)
# instead of silently converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

# Row 2: Bitcoin
Set-TextValue $ws.Range("D2") '63.970.17'

# Row 3: Ethereum
Set-TextValue $ws.Range("D3") '3.466.85'
Set-TextValue $ws.Range("E3") '  -1.07%  '

# Row 4: TetherUSD
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  -0.07%  '

# Row 5: BNB
Set-TextValue $ws.Range("D5") '582.92'
Set-TextValue $ws.Range("E5") '  -0.58%  '

# Row 6: Solana
Set-TextValue $ws.Range("D6") '130.85'
Set-TextValue $ws.Range("E6") '  -2.76%  '

# Row 7: USDC
Set-TextValue $ws.Range("E7") '  +0.00%  '

# Row 8: XRP
Set-TextValue $ws.Range("E8") '  -1.23%  '

# Row 9: Toncoin
Set-TextValue $ws.Range("E9") '  +4.40%  '

# Row 10: Dogecoin
Set-TextValue $ws.Range("E10") '  -2.19%  '

# Row 11: Cardano
Set-TextValue $ws.Range("D11") '0.385'
Set-TextValue $ws.Range("E11") '  -0.40%  '

# Row 12: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") '4.054.31'
Set-TextValue $ws.Range("E12") '  -1.14%  '

# Row 13: TRON
Set-TextValue $ws.Range("D13") '0.119'
Set-TextValue $ws.Range("E13") '  -0.28%  '

# Row 14: WrappedEther
Set-TextValue $ws.Range("B14") 'ShibaInu'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D14") '0.0000177'
Set-TextValue $ws.Range("E14") '  -3.53%  '

# Row 15: ShibaInu
Set-TextValue $ws.Range("B15") 'WrappedEther'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D15") '3.458.48'
Set-TextValue $ws.Range("E15") '  -1.32%  '

# Row 16: WrappedBTC
Set-TextValue $ws.Range("D16") '63.949.85'
Set-TextValue $ws.Range("E16") '  -0.60%  '

# Row 17: Avalanche
Set-TextValue $ws.Range("D17") '24.30'
Set-TextValue $ws.Range("E17") '  -6.56%  '

# Row 18: Uniswap
Set-TextValue $ws.Range("D18") '9.93'
Set-TextValue $ws.Range("E18") '  +0.05%  '

# Row 19: Polkadot
Set-TextValue $ws.Range("D19") '5.67'
Set-TextValue $ws.Range("E19") '  -1.51%  '

# Row 20: Chainlink
Set-TextValue $ws.Range("D20") '13.40'
Set-TextValue $ws.Range("E20") '  -2.18%  '

# Row 21: BitcoinCash
Set-TextValue $ws.Range("D21") '382.64'
Set-TextValue $ws.Range("E21") '  -2.73%  '

# Row 22: Polygon
Set-TextValue $ws.Range("D22") '0.568'
Set-TextValue $ws.Range("E22") '  -0.82%  '

# Row 23: WrappedeETH
Set-TextValue $ws.Range("D23") '3.603.84'
Set-TextValue $ws.Range("E23") '  -1.14%  '

# Row 24: Litecoin
Set-TextValue $ws.Range("D24") '74.52'
Set-TextValue $ws.Range("E24") '  +0.32%  '

# Row 25: Dai
Set-TextValue $ws.Range("E25") '  +0.08%  '

# Row 26: LEO
Set-TextValue $ws.Range("D26") '5.62'

# Row 27: PEPE
Set-TextValue $ws.Range("E27") '  -3.90%  '

# Row 28: Binance-PegBSC-USD
Set-TextValue $ws.Range("D28") '0.999'
Set-TextValue $ws.Range("E28") '  -0.08%  '

# Row 29: PancakeSwap
Set-TextValue $ws.Range("E29") '  -0.58%  '

# Row 30: RenderToken
Set-TextValue $ws.Range("D30") '7.05'
Set-TextValue $ws.Range("E30") '  -4.92%  '

# Row 31: Fetch.AI
Set-TextValue $ws.Range("E31") '  -4.99%  '

# Row 32: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D32") '7.90'
Set-TextValue $ws.Range("E32") '  -4.74%  '

# Row 33: RenzoRestakedETH
Set-TextValue $ws.Range("D33") '3.493.38'
Set-TextValue $ws.Range("E33") '  -0.91%  '

# Row 34: Kaspa
Set-TextValue $ws.Range("D34") '0.152'
Set-TextValue $ws.Range("E34") '  +0.77%  '

# Row 35: USDe
Set-TextValue $ws.Range("E35") '  -0.04%  '

# Row 36: EthereumClassic
Set-TextValue $ws.Range("D36") '22.80'
Set-TextValue $ws.Range("E36") '  -2.85%  '

# Row 37: NEARProtocol
Set-TextValue $ws.Range("D37") '5.17'
Set-TextValue $ws.Range("E37") '  -0.46%  '

# Row 38: Aptos
Set-TextValue $ws.Range("D38") '6.72'
Set-TextValue $ws.Range("E38") '  -3.03%  '

# Row 39: Monero
Set-TextValue $ws.Range("D39") '161.61'
Set-TextValue $ws.Range("E39") '  -1.27%  '

# Row 40: ImmutableX
Set-TextValue $ws.Range("E40") '  -4.78%  '

# Row 41: Hedera
Set-TextValue $ws.Range("D41") '0.0772'
Set-TextValue $ws.Range("E41") '  -1.56%  '

# Row 42: Mantle
Set-TextValue $ws.Range("D42") '0.795'
Set-TextValue $ws.Range("E42") '  -1.54%  '

# Row 43: FirstDigitalUSD
Set-TextValue $ws.Range("E43") '  -0.07%  '

# Row 44: OKB
Set-TextValue $ws.Range("E44") '  -1.09%  '

# Row 45: Filecoin
Set-TextValue $ws.Range("D45") '4.26'
Set-TextValue $ws.Range("E45") '  -3.58%  '

# Row 46: Stacks
Set-TextValue $ws.Range("E46") '  -2.83%  '

# Row 47: EnergySwap
Set-TextValue $ws.Range("D47") '23.20'
Set-TextValue $ws.Range("E47") '  -8.05%  '

# Row 48: ONDO
Set-TextValue $ws.Range("E48") '  -4.86%  '

# Row 49: Cosmos
Set-TextValue $ws.Range("D49") '6.68'
Set-TextValue $ws.Range("E49") '  -1.59%  '

# Row 50: Maker
Set-TextValue $ws.Range("B50") 'SuiNetwork'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D50") '0.895'
Set-TextValue $ws.Range("E50") '  -0.23%  '

# Row 51: SuiNetwork
Set-TextValue $ws.Range("B51") 'Maker'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D51") '2.329.72'
Set-TextValue $ws.Range("E51") '  -5.61%  '
